# Auto-generated PowerShell COM-interop script
# Applies the scheduled-runner price/profit update to each Leve-profit sheet.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3317.8
$ws.Range("I76").Value = 3107.5
$ws.Range("J76").Value = 3738.4
$ws.Range("K76").Value = 3107.5
$ws.Range("L76").Value = 3738.4
$ws.Range("M76").Value = -2792.5
$ws.Range("N76").Value = -4368.4
$ws.Range("H79").Value = 3317.8
$ws.Range("I79").Value = 3107.5
$ws.Range("J79").Value = 3738.4
$ws.Range("K79").Value = 3107.5
$ws.Range("L79").Value = 3738.4
$ws.Range("M79").Value = -2015.5
$ws.Range("N79").Value = -5922.4
$ws.Range("H112").Value = 3745.7273
$ws.Range("I112").Value = 1799
$ws.Range("J112").Value = 4178.3335
$ws.Range("K112").Value = 5397
$ws.Range("L112").Value = 12535.0005
$ws.Range("M112").Value = -4289
$ws.Range("N112").Value = -14751.0005
$ws.Range("H132").Value = 1389.5349
$ws.Range("I132").Value = 1372
$ws.Range("K132").Value = 4116
$ws.Range("M132").Value = -1586
$ws.Range("H138").Value = 39079.926
$ws.Range("J138").Value = 3258.6667
$ws.Range("L138").Value = 9776.000100000001
$ws.Range("N138").Value = -20056.0001
$ws.Range("H140").Value = 86806.55
$ws.Range("J140").Value = 92416.3
$ws.Range("L140").Value = 92416.3
$ws.Range("N140").Value = -102776.3

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4414.095
$ws.Range("I74").Value = 2119.9
$ws.Range("J74").Value = 6499.727
$ws.Range("K74").Value = 2119.9
$ws.Range("L74").Value = 6499.727
$ws.Range("M74").Value = -1245.9
$ws.Range("N74").Value = -8247.726999999999
$ws.Range("H77").Value = 4414.095
$ws.Range("I77").Value = 2119.9
$ws.Range("J77").Value = 6499.727
$ws.Range("K77").Value = 10599.5
$ws.Range("L77").Value = 32498.635
$ws.Range("M77").Value = -6231.5
$ws.Range("N77").Value = -41234.63499999999
$ws.Range("H140").Value = 52914.75
$ws.Range("I140").Value = 375
$ws.Range("J140").Value = 70428
$ws.Range("K140").Value = 375
$ws.Range("L140").Value = 70428
$ws.Range("M140").Value = 4805
$ws.Range("N140").Value = -80788

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 83701.42999999999
$ws.Range("I20").Value = 105620.18
$ws.Range("J20").Value = 3332.6667
$ws.Range("K20").Value = 105620.18
$ws.Range("L20").Value = 3332.6667
$ws.Range("M20").Value = -105373.18
$ws.Range("N20").Value = -3826.6667
$ws.Range("H94").Value = 3207.8
$ws.Range("I94").Value = 3509.75
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 3509.75
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -3058.75
$ws.Range("N94").Value = -2902
$ws.Range("H132").Value = 31383.928
$ws.Range("J132").Value = 31383.928
$ws.Range("L132").Value = 31383.928
$ws.Range("N132").Value = -41503.928
$ws.Range("H135").Value = 103095.2
$ws.Range("J135").Value = 103095.2
$ws.Range("L135").Value = 103095.2
$ws.Range("N135").Value = -113235.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 156.37143
$ws.Range("I7").Value = 53.4
$ws.Range("J7").Value = 233.6
$ws.Range("K7").Value = 53.4
$ws.Range("L7").Value = 233.6
$ws.Range("M7").Value = 59.6
$ws.Range("N7").Value = -459.6
$ws.Range("H31").Value = 2740.3333
$ws.Range("I31").Value = 1475.625
$ws.Range("J31").Value = 4185.7144
$ws.Range("K31").Value = 1475.625
$ws.Range("L31").Value = 4185.7144
$ws.Range("M31").Value = -1180.625
$ws.Range("N31").Value = -4775.7144
$ws.Range("H34").Value = 2740.3333
$ws.Range("I34").Value = 1475.625
$ws.Range("J34").Value = 4185.7144
$ws.Range("K34").Value = 1475.625
$ws.Range("L34").Value = 4185.7144
$ws.Range("M34").Value = -1273.625
$ws.Range("N34").Value = -4589.7144
$ws.Range("H132").Value = 1128968.4
$ws.Range("J132").Value = 1447856.2
$ws.Range("L132").Value = 4343568.6
$ws.Range("N132").Value = -4348628.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 186.33333
$ws.Range("I47").Value = 193.6
$ws.Range("K47").Value = 580.8
$ws.Range("M47").Value = -149.8
$ws.Range("H98").Value = 3205.1177
$ws.Range("I98").Value = 3650.1667
$ws.Range("K98").Value = 10950.5001
$ws.Range("M98").Value = -9452.500100000001
$ws.Range("H122").Value = 918.3570999999999
$ws.Range("J122").Value = 1144.4
$ws.Range("L122").Value = 10299.6
$ws.Range("N122").Value = -15199.6
$ws.Range("H137").Value = 4395.769
$ws.Range("I137").Value = 2405.75
$ws.Range("J137").Value = 7579.8
$ws.Range("K137").Value = 7217.25
$ws.Range("L137").Value = 22739.4
$ws.Range("M137").Value = -2117.25
$ws.Range("N137").Value = -32939.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 254.125
$ws.Range("I2").Value = 202.08333
$ws.Range("J2").Value = 410.25
$ws.Range("K2").Value = 202.08333
$ws.Range("L2").Value = 410.25
$ws.Range("M2").Value = -89.08332999999999
$ws.Range("N2").Value = -636.25
$ws.Range("H14").Value = 125406
$ws.Range("I14").Value = 125406
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 125406
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -125238
$ws.Range("H70").Value = 4270.2856
$ws.Range("I70").Value = 4794.3335
$ws.Range("J70").Value = 3877.25
$ws.Range("K70").Value = 4794.3335
$ws.Range("L70").Value = 3877.25
$ws.Range("M70").Value = -4524.3335
$ws.Range("N70").Value = -4417.25
$ws.Range("H73").Value = 4270.2856
$ws.Range("I73").Value = 4794.3335
$ws.Range("J73").Value = 3877.25
$ws.Range("K73").Value = 4794.3335
$ws.Range("L73").Value = 3877.25
$ws.Range("M73").Value = -3858.3335
$ws.Range("N73").Value = -5749.25
$ws.Range("H97").Value = 1100.8276
$ws.Range("I97").Value = 1103.9
$ws.Range("K97").Value = 1103.9
$ws.Range("M97").Value = -607.9000000000001
$ws.Range("H102").Value = 1558
$ws.Range("I102").Value = 1506.5883
$ws.Range("J102").Value = 1995
$ws.Range("K102").Value = 1506.5883
$ws.Range("L102").Value = 1995
$ws.Range("M102").Value = 115.4117000000001
$ws.Range("N102").Value = -5239
$ws.Range("H132").Value = 5062.222
$ws.Range("I132").Value = 3257.5312
$ws.Range("K132").Value = 9772.5936
$ws.Range("M132").Value = -7242.5936
$ws.Range("H141").Value = 130195
$ws.Range("I141").Value = 110390
$ws.Range("J141").Value = 150000
$ws.Range("K141").Value = 110390
$ws.Range("L141").Value = 150000
$ws.Range("M141").Value = -105210
$ws.Range("N141").Value = -160360
$ws.Range("N14").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3707870.8
$ws.Range("I40").Value = 4464.5713
$ws.Range("K40").Value = 4464.5713
$ws.Range("M40").Value = -4328.5713
$ws.Range("H55").Value = 1834.7727
$ws.Range("J55").Value = 3683.5
$ws.Range("L55").Value = 3683.5
$ws.Range("N55").Value = -4029.5
$ws.Range("H61").Value = 3690.6
$ws.Range("I61").Value = 3900.889
$ws.Range("K61").Value = 3900.889
$ws.Range("M61").Value = -3698.889
$ws.Range("H68").Value = 4206.25
$ws.Range("J68").Value = 3890.6
$ws.Range("L68").Value = 3890.6
$ws.Range("N68").Value = -5388.6
$ws.Range("H71").Value = 4206.25
$ws.Range("J71").Value = 3890.6
$ws.Range("L71").Value = 19453
$ws.Range("N71").Value = -26941
$ws.Range("H113").Value = 3690.6
$ws.Range("I113").Value = 3900.889
$ws.Range("K113").Value = 3900.889
$ws.Range("M113").Value = -1730.889
$ws.Range("H136").Value = 4112.3
$ws.Range("I136").Value = 4376.7856
$ws.Range("K136").Value = 13130.3568
$ws.Range("M136").Value = -10580.3568

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5833
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 6249.5
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 6249.5
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -7497.5
$ws.Range("H65").Value = 5833
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 6249.5
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 31247.5
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -37487.5
$ws.Range("H81").Value = 9933.923000000001
$ws.Range("I81").Value = 11620.5
$ws.Range("J81").Value = 7235.4
$ws.Range("K81").Value = 23241
$ws.Range("L81").Value = 14470.8
$ws.Range("M81").Value = -22180
$ws.Range("N81").Value = -16592.8
$ws.Range("H84").Value = 9933.923000000001
$ws.Range("I84").Value = 11620.5
$ws.Range("J84").Value = 7235.4
$ws.Range("K84").Value = 116205
$ws.Range("L84").Value = 72354
$ws.Range("M84").Value = -110901
$ws.Range("N84").Value = -82962
$ws.Range("H113").Value = 1093.421
$ws.Range("I113").Value = 1314.6666
$ws.Range("K113").Value = 3943.9998
$ws.Range("M113").Value = -1773.9998
$ws.Range("H136").Value = 2549.625
$ws.Range("I136").Value = 2470.5806
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7411.7418
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4861.7418
$ws.Range("N136").Value = -20100
$ws.Range("H140").Value = 149995.5
$ws.Range("J140").Value = 149995.5
$ws.Range("L140").Value = 149995.5
$ws.Range("N140").Value = -160355.5
$ws.Range("H141").Value = 77395.5
$ws.Range("J141").Value = 77395.5
$ws.Range("L141").Value = 77395.5
$ws.Range("N141").Value = -87755.5

Write-Host "Applied scheduled profit/price updates across all sheets."